$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I2").Value = 381.1111
$ws.Range("H2").Value = 400.46667
$ws.Range("J2").Value = 429.5
$ws.Range("L2").Value = 429.5
$ws.Range("N2").Value = -655.5
$ws.Range("M2").Value = -268.1111
$ws.Range("K2").Value = 381.1111
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("M15").Value = -3469.7306
$ws.Range("H15").Value = 1212.9102
$ws.Range("K15").Value = 3638.7306
$ws.Range("I15").Value = 1212.9102
$ws.Range("K29").Value = 949.9999799999999
$ws.Range("I29").Value = 316.66666
$ws.Range("H29").Value = 316.66666
$ws.Range("M29").Value = -668.9999799999999
$ws.Range("H43").Value = 1123.1842
$ws.Range("J43").Value = 1352.1428
$ws.Range("N43").Value = -1490.1428
$ws.Range("L43").Value = 1352.1428
$ws.Range("J53").Value = 1461.7273
$ws.Range("L53").Value = 1461.7273
$ws.Range("N53").Value = -2735.7273
$ws.Range("M53").Value = 104.53845
$ws.Range("K53").Value = 532.46155
$ws.Range("I53").Value = 532.46155
$ws.Range("H53").Value = 958.375
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("I127").Value = 561.5454999999999
$ws.Range("H127").Value = 2454.625
$ws.Range("J127").Value = 2917.3777
$ws.Range("M127").Value = 3275.3635
$ws.Range("L127").Value = 8752.133099999999
$ws.Range("N127").Value = -18672.1331
$ws.Range("K127").Value = 1684.6365
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
$ws.Range("H136").Value = 0
$ws.Range("I137").Value = 2037.1555
$ws.Range("H137").Value = 2044.2122
$ws.Range("J137").Value = 2059.3333
$ws.Range("M137").Value = -3561.4665
$ws.Range("L137").Value = 6177.999899999999
$ws.Range("N137").Value = -11277.9999
$ws.Range("K137").Value = 6111.4665
$ws.Range("K138").Value = 4491.242999999999
$ws.Range("I138").Value = 1497.081
$ws.Range("J138").Value = 9350
$ws.Range("H138").Value = 4067.1272
$ws.Range("M138").Value = 648.7570000000005
$ws.Range("L138").Value = 28050
$ws.Range("N138").Value = -38330

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I25").Value = 1550
$ws.Range("H25").Value = 1917
$ws.Range("J25").Value = 3018
$ws.Range("M25").Value = -1148
$ws.Range("L25").Value = 3018
$ws.Range("N25").Value = -3822
$ws.Range("K25").Value = 1550
$ws.Range("I32").Value = 2967.6511
$ws.Range("H32").Value = 3336815.8
$ws.Range("J32").Value = 23816168
$ws.Range("M32").Value = -2680.6511
$ws.Range("L32").Value = 23816168
$ws.Range("N32").Value = -23816742
$ws.Range("K32").Value = 2967.6511
$ws.Range("I33").Value = 33000
$ws.Range("H33").Value = 33000
$ws.Range("M33").Value = -32671
$ws.Range("K33").Value = 33000
$ws.Range("M37").Value = -99721
$ws.Range("K37").Value = 99994
$ws.Range("I37").Value = 99994
$ws.Range("H37").Value = 41331.332
$ws.Range("N61").Value = -2407.375
$ws.Range("M61").Value = -36290864
$ws.Range("L61").Value = 1983.375
$ws.Range("K61").Value = 36291076
$ws.Range("I61").Value = 36291076
$ws.Range("J61").Value = 1983.375
$ws.Range("H61").Value = 23937342
$ws.Range("N74").Value = -3889.1177
$ws.Range("M74").Value = 63.17645000000005
$ws.Range("L74").Value = 2141.1177
$ws.Range("K74").Value = 810.82355
$ws.Range("H74").Value = 1254.2549
$ws.Range("I74").Value = 810.82355
$ws.Range("J74").Value = 2141.1177
$ws.Range("K77").Value = 4054.11775
$ws.Range("I77").Value = 810.82355
$ws.Range("H77").Value = 1254.2549
$ws.Range("J77").Value = 2141.1177
$ws.Range("M77").Value = 313.8822500000001
$ws.Range("L77").Value = 10705.5885
$ws.Range("N77").Value = -19441.5885
$ws.Range("J133").Value = 40000
$ws.Range("H133").Value = 40000
$ws.Range("L133").Value = 40000
$ws.Range("N133").Value = -45060
$ws.Range("J136").Value = 1983.375
$ws.Range("M136").Value = -108870678
$ws.Range("L136").Value = 5950.125
$ws.Range("N136").Value = -11050.125
$ws.Range("K136").Value = 108873228
$ws.Range("I136").Value = 36291076
$ws.Range("H136").Value = 23937342

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N24").Value = -1470
$ws.Range("M24").Value = -173
$ws.Range("K24").Value = 408
$ws.Range("I24").Value = 408
$ws.Range("H24").Value = 704
$ws.Range("J24").Value = 1000
$ws.Range("L24").Value = 1000
$ws.Range("L36").Value = 18002
$ws.Range("N36").Value = -19070
$ws.Range("M36").Value = -125.25
$ws.Range("K36").Value = 659.25
$ws.Range("I36").Value = 659.25
$ws.Range("H36").Value = 4127.8
$ws.Range("J36").Value = 18002
$ws.Range("J37").Value = 1150
$ws.Range("L37").Value = 1150
$ws.Range("N37").Value = -1424
$ws.Range("M37").Value = -671.7143
$ws.Range("K37").Value = 808.7143
$ws.Range("I37").Value = 808.7143
$ws.Range("H37").Value = 911.1
$ws.Range("I107").Value = 27778670
$ws.Range("J107").Value = 3103.25
$ws.Range("H107").Value = 19232342
$ws.Range("N107").Value = -6943.25
$ws.Range("M107").Value = -27776750
$ws.Range("L107").Value = 3103.25
$ws.Range("K107").Value = 27778670
$ws.Range("L140").Value = 48944.445
$ws.Range("N140").Value = -59304.445
$ws.Range("H140").Value = 48944.445
$ws.Range("J140").Value = 48944.445

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I31").Value = 856.8570999999999
$ws.Range("J31").Value = 1302.1724
$ws.Range("H31").Value = 1115.14
$ws.Range("N31").Value = -1892.1724
$ws.Range("M31").Value = -561.8570999999999
$ws.Range("L31").Value = 1302.1724
$ws.Range("K31").Value = 856.8570999999999
$ws.Range("K34").Value = 856.8570999999999
$ws.Range("I34").Value = 856.8570999999999
$ws.Range("H34").Value = 1115.14
$ws.Range("J34").Value = 1302.1724
$ws.Range("L34").Value = 1302.1724
$ws.Range("N34").Value = -1706.1724
$ws.Range("M34").Value = -654.8570999999999
$ws.Range("K134").Value = 2706.7242
$ws.Range("I134").Value = 902.2414
$ws.Range("H134").Value = 14286882
$ws.Range("J134").Value = 83335784
$ws.Range("M134").Value = -171.7242000000001
$ws.Range("L134").Value = 250007352
$ws.Range("N134").Value = -250012422

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("K5").Value = 205883550
$ws.Range("I5").Value = 68627850
$ws.Range("H5").Value = 37638220
$ws.Range("J5").Value = 7955.357
$ws.Range("N5").Value = -24090.071
$ws.Range("M5").Value = -205883438
$ws.Range("L5").Value = 23866.071
$ws.Range("L68").Value = 29697.999
$ws.Range("N68").Value = -31319.999
$ws.Range("M68").Value = -857.2904000000001
$ws.Range("K68").Value = 1668.2904
$ws.Range("I68").Value = 556.0968
$ws.Range("H68").Value = 3163.5117
$ws.Range("J68").Value = 9899.333000000001
$ws.Range("K71").Value = 5004.8712
$ws.Range("I71").Value = 556.0968
$ws.Range("H71").Value = 3163.5117
$ws.Range("J71").Value = 9899.333000000001
$ws.Range("L71").Value = 89093.997
$ws.Range("N71").Value = -97205.997
$ws.Range("M71").Value = -948.8712000000005
$ws.Range("I131").Value = 524.9
$ws.Range("H131").Value = 820.54
$ws.Range("J131").Value = 853.3889
$ws.Range("L131").Value = 2560.1667
$ws.Range("N131").Value = -12640.1667
$ws.Range("M131").Value = 3465.3
$ws.Range("K131").Value = 1574.7
$ws.Range("I135").Value = 68627850
$ws.Range("H135").Value = 37638220
$ws.Range("J135").Value = 7955.357
$ws.Range("M135").Value = -617648115
$ws.Range("L135").Value = 71598.213
$ws.Range("N135").Value = -76668.213
$ws.Range("K135").Value = 617650650

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L132").Value = 34682.769
$ws.Range("K132").Value = 8265.1032
$ws.Range("I132").Value = 2755.0344
$ws.Range("H132").Value = 5480.6665
$ws.Range("J132").Value = 11560.923
$ws.Range("N132").Value = -39742.769
$ws.Range("M132").Value = -5735.1032

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M16").Value = -506
$ws.Range("H16").Value = 676
$ws.Range("K16").Value = 676
$ws.Range("I16").Value = 676
$ws.Range("K34").Value = 7999
$ws.Range("I34").Value = 7999
$ws.Range("H34").Value = 8999.5
$ws.Range("M34").Value = -7827
$ws.Range("M61").Value = -1798
$ws.Range("K61").Value = 2000
$ws.Range("I61").Value = 2000
$ws.Range("H61").Value = 2000
$ws.Range("K93").Value = 1011.38464
$ws.Range("I93").Value = 1011.38464
$ws.Range("H93").Value = 1025.88
$ws.Range("J93").Value = 1041.5834
$ws.Range("M93").Value = 236.61536
$ws.Range("L93").Value = 1041.5834
$ws.Range("N93").Value = -3537.5834
$ws.Range("M113").Value = 170
$ws.Range("K113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("H113").Value = 2000
$ws.Range("J136").Value = 23811626
$ws.Range("M136").Value = -170336580
$ws.Range("L136").Value = 71434878
$ws.Range("N136").Value = -71439978
$ws.Range("K136").Value = 170339130
$ws.Range("I136").Value = 56779710
$ws.Range("H136").Value = 36417070

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J31").Value = 3660
$ws.Range("H31").Value = 3660
$ws.Range("N31").Value = -4356
$ws.Range("L31").Value = 3660
$ws.Range("K34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("M34").ClearContents()
$ws.Range("M113").Value = -200008346
$ws.Range("K113").Value = 200010516
$ws.Range("I113").Value = 66670172
$ws.Range("H113").Value = 58826828
$ws.Range("J113").Value = 1750
$ws.Range("L113").Value = 5250
$ws.Range("N113").Value = -9590
$ws.Range("J136").Value = 2850
$ws.Range("M136").Value = -65218140
$ws.Range("L136").Value = 8550
$ws.Range("N136").Value = -13650
$ws.Range("K136").Value = 65220690
$ws.Range("I136").Value = 21740230
$ws.Range("H136").Value = 11629821
